$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 blog cards cycle: "ser: 95" drops out, "ser: 96" -> I7 becomes "ser: 96",
# "ser: 98" shifts from E7 to ... etc. Net effect on the visible text values:
#   C7: ser 98 -> ser 97  (new day's blog post)
#   E7: ser 96 -> ser 98
#   I7: ser 95 -> ser 96
$ws.Range("C7").Value = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 97"
$ws.Range("E7").Value = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 98"
$ws.Range("I7").Value = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 96"
